# Append latest Lancers scrape run (2025-09-24 12:36 JST) to the "ランサーズ"
# sheet: refresh the "取得日時" timestamp on every existing row, fold in the
# newly-scraped postings (re-sorted by priority score), and widen column D to
# fit the new price strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D ("価格") needs to go from a raw OOXML width of 28 to 32. The
# ColumnWidth setter works in "Excel character" units and rounds to the
# nearest pixel, so 31.16 is the value that lands exactly on a stored
# width of 32 (mirrors how 28 previously corresponded to ~27.17).
$ws.Columns.Item(4).ColumnWidth = 31.16

$timestamp = "2025-09-24 12:36:20"
$category = "システム開発"
$deadline = "期限情報なし"

# Every F-column cell's URL changes (rows shift to a new score-sorted order),
# so drop all existing hyperlinks up front and re-add them fresh below —
# mutating a hyperlink's .Address in place (or re-Add-ing over an existing
# one) leaves the old relationship dangling and duplicates the entry instead
# of replacing it.
$ws.Hyperlinks.Delete()

# r, title, price, url, score, skills (skills = $null when the column is blank)
$rows = @(
    @(2,  "【AI活用】データ分析Webサービス開発パートナー募集",                                              "200,000 円 ~ 300,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399092", 368, "🔥AI,Ai ◆開発"),
    @(3,  "あなたAIクローン構築パートナー募集・モデル制作&新規依頼",                                        "100,000 円 ~ 200,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399534", 303, "🔥AI,Ai"),
    @(4,  "Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集",                  "200,000 円 ~ 300,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399398", 298, "🔥Python ◆開発,スクレイピング"),
    @(5,  "【急募】カスタマー向けFAQチャットbotの開発依頼",                                                  "1,000,000 円 ~ 3,000,000 円 / 固定", "https://www.lancers.jp/work/detail/5399558", 180, "★bot ◆開発"),
    @(6,  "既存Excelをベースにした短期計画書管理のWebシステム開発",                                          "100,000 円 ~ 200,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399602", 153, "◆開発,システム開発 ◇管理"),
    @(7,  "【RPA/ブラウザ操作自動化】Webフォーム大量登録の自動化(継続依頼あり)",                              "200,000 円 ~ 300,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399631", 88,  "◆自動化"),
    @(8,  "Googleフォーム × スプレッドシート × GAS 自動化(ストレスチェック診断/台帳保存あり)",               "50,000 円 ~ 100,000 円 / 固定",      "https://www.lancers.jp/work/detail/5399200", 88,  "◆自動化"),
    @(9,  "急募 【急募】Excelで株の保有リストを自動化したいので制作してくださる方募集!",                      "20,000 円 ~ 50,000 円 / 固定",       "https://www.lancers.jp/work/detail/5399727", 83,  "◆自動化"),
    @(10, "【急募】住宅展示場マッチング診断サービスのMVP開発依頼",                                           "500,000 円 ~ 1,000,000 円 / 固定",   "https://www.lancers.jp/work/detail/5399759", 75,  "◆開発"),
    @(11, "完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします",       "~ 5,000 円 / 固定",                  "https://www.lancers.jp/work/detail/5399071", 70,  "◆効率化"),
    @(12, "【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼",                                           "200,000 円 ~ 300,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399313", 68,  "◆開発"),
    @(13, "〖リモート可〗Delphiエンジニア募集",                                                              "300,000 円 ~ 500,000 円 / 固定",     "https://www.lancers.jp/work/detail/5341051", 25,  $null),
    @(14, "【急募】フロントエンドエンジニア募集!(ややWEBコーダー寄り)",                                     "300,000 円 ~ 500,000 円 / 固定",     "https://www.lancers.jp/work/detail/5399545", 25,  $null),
    @(15, "【相談から実装まで伴走できる方歓迎】介護・福祉×テクノロジー事例収集の仕組みづくり",               "50,000 円 ~ 100,000 円 / 固定",      "https://www.lancers.jp/work/detail/5398932", 18,  $null),
    @(16, "限定公開 PR 限定公開の仕事",                                                                      "20,000 円 ~ 50,000 円 / 固定",       "https://www.lancers.jp/work/detail/5399347", 13,  $null),
    @(17, "Android kotlin 非同期処理の呼び方",                                                               "5,000 円 ~ 10,000 円 / 固定",        "https://www.lancers.jp/work/detail/5399765", 10,  $null)
)

foreach ($row in $rows) {
    $r      = $row[0]
    $title  = $row[1]
    $price  = $row[2]
    $url    = $row[3]
    $score  = $row[4]
    $skills = $row[5]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $title
    $ws.Cells.Item($r, 3).Value = $category
    $ws.Cells.Item($r, 4).Value = $price
    $ws.Cells.Item($r, 5).Value = $deadline
    $ws.Cells.Item($r, 6).Value = $url
    $ws.Cells.Item($r, 7).Value = $score

    if ($skills) {
        $ws.Cells.Item($r, 8).Value = $skills
    }

    # Re-applying the built-in "Hyperlink" named style after Add keeps every
    # F cell on the original style index instead of minting a duplicate
    # (functionally-identical) style per call.
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url)
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
